$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Artista"
$ws.Range("B1").Value = "Diva"
$ws.Range("C1").Value = "OA"
$ws.Range("D1").Value = "Wildcard"
$ws.Range("E1").Value = "Achiever"
$ws.Range("F1").Value = "EMO"
$ws.Range("G1").Value = "Gamer"
$ws.Range("H1").Value = "Softie"
